# Apply the "2D loop" presence update to the Algoritmi attendance sheet.
# Column G corresponds to "saapt. 6" (week 6) attendance marks; several
# students now have a "2" recorded instead of "1", and one student (row 9)
# gets a new mark of "1" in week 6 that previously had no entry at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update week 6 (column G) attendance values.
$ws.Range("G3").Value  = 2
$ws.Range("G6").Value  = 2
$ws.Range("G7").Value  = 2
$ws.Range("G9").Value  = 1
$ws.Range("G12").Value = 2
$ws.Range("G13").Value = 2
$ws.Range("G14").Value = 2
$ws.Range("G17").Value = 2
$ws.Range("G18").Value = 2
$ws.Range("G19").Value = 2
$ws.Range("G21").Value = 2

# Recalculate so the Q-column "Prezente" totals (SUM formulas) refresh.
$excel.Calculate()

# Move the active selection on the frozen pane to I8 (was I12).
$ws.Range("I8").Select()
